$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.730.74'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.631.26'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.00'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.09%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.56'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('E10').Value = '  +1.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.376'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.33%  '
$ws.Range('E12').Value = '  +4.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.098.27'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.26'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +12.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.713.24'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000143'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.634.26'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('E18').Value = '  +2.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.73'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '349.25'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.90'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.533'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('E26').Value = '  +2.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.75%  '
$ws.Range('E28').Value = '  +11.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0803'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.66'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '169.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.20%  '
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.64'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.07'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.45'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.32%  '
$ws.Range('E36').Value = '  +8.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.65'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '332.90'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +13.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.01'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '38.40'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.86%  '
$ws.Range('E41').Value = '  +3.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.16'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.99'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.84%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '133.67'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.97%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0998'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '20.12'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.07%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0558'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.41%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.613'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.81%  '
$ws.Range('E50').Value = '  +1.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.60%  '
